{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// --- 1 & 2) Title + author placeholders -----------------------------------\n// \"La Palma Earthquakes\"  -> \"Test of Quarto Manuscript\"   (style \"Title\")\n// \"Steve Purves\"          -> \"Author 1\"                    (1st \"Author\" para)\n// \"Rowan Cockett\"         -> \"Author 2\"                    (2nd \"Author\" para)\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"style\");\n}\nawait context.sync();\n\nlet authorCount = 0;\nfor (const p of paragraphs.items) {\n  const style = (p.style || \"\").trim();\n\n  if (style === \"Title\") {\n    p.insertText(\"Test of Quarto Manuscript\", \"Replace\");\n  } else if (style === \"Author\") {\n    authorCount += 1;\n    p.insertText(`Author ${authorCount}`, \"Replace\");\n  }\n}\nawait context.sync();\n\n// --- 3) Collapse the long sentence after \"Marrero et al. (2019)\" ----------\nconst oldSentence =\n  \", have proposed that there are two main magma reservoirs feeding the Cumbre Vieja volcano; one in the mantle (30-40km depth) which charges and in turn feeds a shallower crustal reservoir (10-20km depth).\";\n\nconst found = body.search(oldSentence, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(\"\u2026\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $doc / $app resolve to the live session; the document is open\n# as $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1        # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Title: \"La Palma Earthquakes\" -> \"Test of Quarto Manuscript\"\nReplace-Text \"La Palma Earthquakes\" \"Test of Quarto Manuscript\"\n\n# 2) Authors: \"Steve Purves\" -> \"Author 1\", \"Rowan Cockett\" -> \"Author 2\"\nReplace-Text \"Steve Purves\" \"Author 1\"\nReplace-Text \"Rowan Cockett\" \"Author 2\"\n\n# 3) Collapse the long sentence after \"Marrero et al. (2019)\" down to an ellipsis.\nReplace-Text \", have proposed that there are two main magma reservoirs feeding the Cumbre Vieja volcano; one in the mantle (30-40km depth) which charges and in turn feeds a shallower crustal reservoir (10-20km depth).\" \"\u2026\"\n"}
